$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column C ---
# This shifts the old C:I (Request Type .. PM) one column to the right (D:J)
# and pushes the used range from A1:I4 to A1:J4.
$ws.Columns.Item(3).Insert()

# --- Column widths (A:B share the old "A" width, C gets the old "B" width,
#     D:I inherit the old C:H widths automatically via the insert/shift above,
#     but explicitly restate them to be safe) ---
$ws.Columns.Item(1).ColumnWidth = 13.944010416666666
$ws.Columns.Item(2).ColumnWidth = 13.944010416666666
$ws.Columns.Item(3).ColumnWidth = 23.276041666666668
$ws.Columns.Item(4).ColumnWidth = 14.944010416666666
$ws.Columns.Item(5).ColumnWidth = 17.944010416666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668
$ws.Columns.Item(7).ColumnWidth = 17.944010416666668
$ws.Columns.Item(8).ColumnWidth = 15.498697916666666
$ws.Columns.Item(9).ColumnWidth = 15.944010416666666

# --- New header cell C3 = "EmailAddress" ---
# Reset to Normal first so it doesn't inherit column B's numeric format.
$ws.Range("C3").Style = "Normal"
$ws.Range("C3").Value = "EmailAddress"
$ws.Range("C3").Font.Bold = $true
$ws.Range("C3").Font.Size = 12
$ws.Range("C3").Font.Name = "Calibri"
$ws.Range("C3").HorizontalAlignment = -4108

# --- Selection, matching the saved workbook's cursor position ---
$ws.Range("C11").Select() | Out-Null
